$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.996.77'
$ws.Range('E2').Value = '  +0.58%  '
$ws.Range('D3').Value = '2.662.06'
$ws.Range('E3').Value = '  +1.86%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '529.98'
$ws.Range('E5').Value = '  +3.13%  '
$ws.Range('D6').Value = '155.66'
$ws.Range('E6').Value = '  +0.70%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '0.582'
$ws.Range('E8').Value = '  -1.27%  '
$ws.Range('D9').Value = '6.49'
$ws.Range('E9').Value = '  -3.27%  '
$ws.Range('E10').Value = '  +4.76%  '
$ws.Range('E11').Value = '  +1.24%  '
$ws.Range('E12').Value = '  -0.44%  '
$ws.Range('D13').Value = '3.122.19'
$ws.Range('E13').Value = '  +1.65%  '
$ws.Range('D14').Value = '60.994.41'
$ws.Range('E14').Value = '  +0.64%  '
$ws.Range('D15').Value = '22.12'
$ws.Range('E15').Value = '  +1.91%  '
$ws.Range('E16').Value = '  +1.48%  '
$ws.Range('D17').Value = '2.671.27'
$ws.Range('E17').Value = '  +1.86%  '
$ws.Range('D18').Value = '4.78'
$ws.Range('E18').Value = '  +0.75%  '
$ws.Range('D19').Value = "'354.30"
$ws.Range('E19').Value = '  -1.00%  '
$ws.Range('D20').Value = '10.67'
$ws.Range('E20').Value = '  +0.33%  '
$ws.Range('D21').Value = '6.33'
$ws.Range('D22').Value = "'1.00"
$ws.Range('E22').Value = '  +0.17%  '
$ws.Range('D23').Value = '61.67'
$ws.Range('E23').Value = '  +1.49%  '
$ws.Range('E24').Value = '  +1.44%  '
$ws.Range('E25').Value = '  +0.73%  '
$ws.Range('E26').Value = '  +0.13%  '
$ws.Range('E27').Value = '  +1.27%  '
$ws.Range('D28').Value = '7.32'
$ws.Range('E28').Value = '  -0.74%  '
$ws.Range('E29').Value = '  -0.02%  '
$ws.Range('E30').Value = '  +3.50%  '
$ws.Range('D31').Value = '19.53'
$ws.Range('E31').Value = '  +0.25%  '
$ws.Range('E32').Value = '  +2.46%  '
$ws.Range('D33').Value = '150.11'
$ws.Range('E33').Value = '  -1.23%  '
$ws.Range('D34').Value = '4.12'
$ws.Range('E34').Value = '  +2.81%  '
$ws.Range('E35').Value = '  +0.21%  '
$ws.Range('D36').Value = '0.919'
$ws.Range('E36').Value = '  +7.97%  '
$ws.Range('D37').Value = '0.895'
$ws.Range('E37').Value = '  +1.28%  '
$ws.Range('D38').Value = '36.89'
$ws.Range('E38').Value = '  +1.47%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').Value = '1.49'
$ws.Range('E39').Value = '  +0.03%  '
$ws.Range('B40').Value = 'Bittensor'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D40').Value = '305.57'
$ws.Range('E40').Value = '  +4.13%  '
$ws.Range('E41').Value = '  +0.54%  '
$ws.Range('E42').Value = '  +3.47%  '
$ws.Range('E43').Value = '  +0.40%  '
$ws.Range('D44').Value = '20.43'
$ws.Range('E44').Value = '  +2.96%  '
$ws.Range('D45').Value = '0.0564'
$ws.Range('E45').Value = '  +1.18%  '
$ws.Range('E46').Value = '  +0.07%  '
$ws.Range('D47').Value = '0.0242'
$ws.Range('E47').Value = '  +2.79%  '
$ws.Range('D48').Value = '4.89'
$ws.Range('E48').Value = '  -1.06%  '
$ws.Range('E49').Value = '  +6.02%  '
$ws.Range('E50').Value = '  +0.72%  '
$ws.Range('D51').Value = '1.999.20'
$ws.Range('E51').Value = '  +0.17%  '
